$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new data rows right before the current row 1010. This shifts the
# existing rows 1010-1097 down to 1013-1100 (matching the dimension growth
# from A1:R1097 to A1:R1100 in the target workbook).
$ws.Rows("1010:1012").Insert()

# New row 1010: Asterix, 1a (guarda), Región de O'Higgins
$ws.Cells.Item(1010, 1).Value = 3
$ws.Cells.Item(1010, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1010, 3).Value = "Coquimbo"
$ws.Cells.Item(1010, 4).Value = 45106
$ws.Cells.Item(1010, 5).Value = 5
$ws.Cells.Item(1010, 6).Value = 100114001
$ws.Cells.Item(1010, 7).Value = "Papa"
$ws.Cells.Item(1010, 8).Value = "Asterix"
$ws.Cells.Item(1010, 9).Value = "1a (guarda)"
$ws.Cells.Item(1010, 10).Value = 480
$ws.Cells.Item(1010, 11).Value = 17000
$ws.Cells.Item(1010, 12).Value = 18000
$ws.Cells.Item(1010, 13).Value = 17479
$ws.Cells.Item(1010, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(1010, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(1010, 16).Value = 699
$ws.Cells.Item(1010, 17).Value = 25
$ws.Cells.Item(1010, 18).Value = "Hortaliza"

# New row 1011: Rosara, 1a (guarda), Provincia de Santiago
$ws.Cells.Item(1011, 1).Value = 3
$ws.Cells.Item(1011, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1011, 3).Value = "Coquimbo"
$ws.Cells.Item(1011, 4).Value = 45106
$ws.Cells.Item(1011, 5).Value = 5
$ws.Cells.Item(1011, 6).Value = 100114001
$ws.Cells.Item(1011, 7).Value = "Papa"
$ws.Cells.Item(1011, 8).Value = "Rosara"
$ws.Cells.Item(1011, 9).Value = "1a (guarda)"
$ws.Cells.Item(1011, 10).Value = 450
$ws.Cells.Item(1011, 11).Value = 17000
$ws.Cells.Item(1011, 12).Value = 18000
$ws.Cells.Item(1011, 13).Value = 17489
$ws.Cells.Item(1011, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(1011, 15).Value = "Provincia de Santiago"
$ws.Cells.Item(1011, 16).Value = 700
$ws.Cells.Item(1011, 17).Value = 25
$ws.Cells.Item(1011, 18).Value = "Hortaliza"

# New row 1012: Rosara, 1a nueva(o), Provincia de Quillota
$ws.Cells.Item(1012, 1).Value = 3
$ws.Cells.Item(1012, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1012, 3).Value = "Coquimbo"
$ws.Cells.Item(1012, 4).Value = 45106
$ws.Cells.Item(1012, 5).Value = 5
$ws.Cells.Item(1012, 6).Value = 100114001
$ws.Cells.Item(1012, 7).Value = "Papa"
$ws.Cells.Item(1012, 8).Value = "Rosara"
$ws.Cells.Item(1012, 9).Value = "1a nueva(o)"
$ws.Cells.Item(1012, 10).Value = 340
$ws.Cells.Item(1012, 11).Value = 17500
$ws.Cells.Item(1012, 12).Value = 18000
$ws.Cells.Item(1012, 13).Value = 17765
$ws.Cells.Item(1012, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(1012, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(1012, 16).Value = 711
$ws.Cells.Item(1012, 17).Value = 25
$ws.Cells.Item(1012, 18).Value = "Hortaliza"
